# Update "想去人数" (F column) figures on sheet "展览" and sheet "全部类型".
# Both sheets track the same events; the "全部类型" rows mirror the
# "展览" rows at different row offsets.

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsAll     = $wb.Worksheets.Item("全部类型")

# 展览 sheet updates
$wsExhibit.Range("F5").Value  = 21
$wsExhibit.Range("F7").Value  = 294
$wsExhibit.Range("F9").Value  = 1032
$wsExhibit.Range("F14").Value = 13219
$wsExhibit.Range("F18").Value = 5435
$wsExhibit.Range("F19").Value = 5562
$wsExhibit.Range("F20").Value = 30

# 全部类型 sheet updates (same events, different rows)
$wsAll.Range("F12").Value = 21
$wsAll.Range("F23").Value = 294
$wsAll.Range("F31").Value = 1032
$wsAll.Range("F36").Value = 13219
$wsAll.Range("F41").Value = 5435
$wsAll.Range("F42").Value = 5562
$wsAll.Range("F43").Value = 30
